$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B251").Value = 187
$ws.Range("C251").Value = 628

$dates = @("08-09-2021","09-09-2021","10-09-2021","11-09-2021","12-09-2021","13-09-2021","14-09-2021")
$row = 252
foreach ($d in $dates) {
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $d
    $ws.Cells.Item($row, 2).Value = 187
    $ws.Cells.Item($row, 3).Value = 628
    $ws.Cells.Item($row, 4).Value = 3940
    $ws.Cells.Item($row, 5).Value = 30
    $row++
}

$ws.Cells.Item(259, 1).NumberFormat = "@"
$ws.Cells.Item(259, 1).Value = "15-09-2021"
$ws.Cells.Item(259, 4).Value = 3940
$ws.Cells.Item(259, 5).Value = 30
